$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.468.31'
$ws.Range("D2").ClearFormats()

$ws.Range("E2").Value = '  -0.88%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.618.48'
$ws.Range("D3").ClearFormats()

$ws.Range("E3").Value = '  -1.72%  '

$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.00'
$ws.Range("D5").ClearFormats()

$ws.Range("E5").Value = '  -1.06%  '

$ws.Range("E6").Value = '  -1.36%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.78'
$ws.Range("D8").ClearFormats()

$ws.Range("E8").Value = '  -1.28%  '

$ws.Range("E9").Value = '  +2.07%  '

$ws.Range("E10").Value = '  -0.05%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0885'
$ws.Range("D11").ClearFormats()

$ws.Range("E11").Value = '  -0.41%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.848.28'
$ws.Range("D12").ClearFormats()

$ws.Range("E12").Value = '  -1.68%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.616.77'
$ws.Range("D13").ClearFormats()

$ws.Range("E13").Value = '  -1.75%  '

$ws.Range("E14").Value = '  -0.42%  '

$ws.Range("E15").Value = '  -2.45%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.17'
$ws.Range("D16").ClearFormats()

$ws.Range("E16").Value = '  +1.48%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.467.74'
$ws.Range("D17").ClearFormats()

$ws.Range("E17").Value = '  -0.75%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '230.88'
$ws.Range("D18").ClearFormats()

$ws.Range("E18").Value = '  +0.23%  '

$ws.Range("E19").Value = '  -0.88%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.51'
$ws.Range("D20").ClearFormats()

$ws.Range("E20").Value = '  -1.71%  '

$ws.Range("E21").Value = '  +0.15%  '

$ws.Range("E22").Value = '  -0.98%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.17'
$ws.Range("D23").ClearFormats()

$ws.Range("E23").Value = '  +1.51%  '

$ws.Range("E24").Value = '  +6.00%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.38'
$ws.Range("D25").ClearFormats()

$ws.Range("E25").Value = '  +0.89%  '

$ws.Range("B26").Value = 'Cosmos'

$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.84'
$ws.Range("D26").ClearFormats()

$ws.Range("E26").Value = '  -1.84%  '

$ws.Range("B27").Value = 'Stellar'

$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.111'
$ws.Range("D27").ClearFormats()

$ws.Range("E27").Value = '  -1.06%  '

$ws.Range("E28").Value = '  +0.21%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.57'
$ws.Range("D29").ClearFormats()

$ws.Range("E29").Value = '  -0.51%  '

$ws.Range("E30").Value = '  -0.75%  '

$ws.Range("E31").Value = '  -0.17%  '

$ws.Range("E32").Value = '  -1.15%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.449.94'
$ws.Range("D33").ClearFormats()

$ws.Range("E33").Value = '  +0.64%  '

$ws.Range("E34").Value = '  -3.33%  '

$ws.Range("E35").Value = '  -3.37%  '

$ws.Range("E36").Value = '  -0.19%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.940'
$ws.Range("D37").ClearFormats()

$ws.Range("E37").Value = '  +3.74%  '

$ws.Range("B38").Value = 'VeChain'

$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0168'
$ws.Range("D38").ClearFormats()

$ws.Range("E38").Value = '  -0.30%  '

$ws.Range("B39").Value = 'ImmutableX'

$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.560'
$ws.Range("D39").ClearFormats()

$ws.Range("E39").Value = '  -1.97%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.863'
$ws.Range("D40").ClearFormats()

$ws.Range("E40").Value = '  -2.33%  '

$ws.Range("E41").Value = '  +0.15%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '67.85'
$ws.Range("D42").ClearFormats()

$ws.Range("E42").Value = '  +3.38%  '

$ws.Range("B43").Value = 'mCoin'

$ws.Range("C43").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.49'
$ws.Range("D43").ClearFormats()

$ws.Range("E43").Value = '  +0.70%  '

$ws.Range("B44").Value = 'WEMIXToken'

$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.991'
$ws.Range("D44").ClearFormats()

$ws.Range("E44").Value = '  -4.03%  '

$ws.Range("B45").Value = 'FraxShare'

$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.41'
$ws.Range("D45").ClearFormats()

$ws.Range("E45").Value = '  -4.61%  '

$ws.Range("B46").Value = 'MXToken'

$ws.Range("C46").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.21'
$ws.Range("D46").ClearFormats()

$ws.Range("E46").Value = '  -2.27%  '

$ws.Range("B47").Value = 'RocketPoolETH'

$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.758.41'
$ws.Range("D47").ClearFormats()

$ws.Range("E47").Value = '  -1.69%  '

$ws.Range("B48").Value = 'RenderToken'

$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.70'
$ws.Range("D48").ClearFormats()

$ws.Range("E48").Value = '  +0.94%  '

$ws.Range("B49").Value = 'Quant'

$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '86.56'
$ws.Range("D49").ClearFormats()

$ws.Range("E49").Value = '  +0.19%  '

$ws.Range("B50").Value = 'BabyDogeCoin'

$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0106'
$ws.Range("D50").ClearFormats()

$ws.Range("E50").Value = '  +17.65%  '

$ws.Range("B51").Value = 'Algorand'

$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.101'
$ws.Range("D51").ClearFormats()

$ws.Range("E51").Value = '  +1.73%  '
